# Swap the 2nd step ("#2") content between TC2 and TC3.
# Before:
#   TC2 step2: B20 = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
#              D20 = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá
#                      constar o nome do usuário logado (que se atribuiu como responsável pela
#                      liquidação) no campo de atribuição (no caso de desatribuição, o nome
#                      deverá ser removido)."
#   TC3 step2: B28 = "Chefe Clica para realizar a liquidação."
#              D28 = "SYSTEM Apresenta a tela de Registrar Liquidações"
#
# After (per diff): the two steps are swapped, so TC2's 2nd step becomes the
# "realizar a liquidação" step, and TC3's 2nd step becomes the
# "atribuir/desatribuir" step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$attribText = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$attribResultText = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$liquidacaoText = "Chefe Clica para realizar a liquidação."
$liquidacaoResultText = "SYSTEM Apresenta a tela de Registrar Liquidações"

# TC2's second step (row 20) now holds the "realizar a liquidação" content.
$ws.Range("B20").Value = $liquidacaoText
$ws.Range("D20").Value = $liquidacaoResultText

# TC3's second step (row 28) now holds the "atribuir/desatribuir" content.
$ws.Range("B28").Value = $attribText
$ws.Range("D28").Value = $attribResultText
